$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) text values.
# Values are stored as text in the original workbook (t="inlineStr"),
# so NumberFormat is forced to "@" (Text) before assignment to avoid
# Excel auto-converting the numeric-looking strings into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "291.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.97%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.73%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.885"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.35%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07262"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.81%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.349"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "28.84%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.686"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.24%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.708"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.25%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8965"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.29%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1662"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.38%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07918"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.19%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08041"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.86%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03099"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.67%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1004"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.30%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001496"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.05%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005848"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.43%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.472"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.36%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.16%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3320"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.50%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1299"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.48%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.994"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-8.65%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2100"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.87%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04505"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.55%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001208"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.31%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004413"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "9.25%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001299"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.74%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01569"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.31%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04385"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.52%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007345"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.88%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009887"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1311"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.95%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002078"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.98%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009348"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-16.52%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005754"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.21%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.13%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "5.41%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002899"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "19.04%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.13%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.13%"
